$wb = $excel.ActiveWorkbook

# Sheet "Overview": update the generate-date string (col G) for rows 7,10-14
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in 7,10,11,12,13,14) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-19 12:20:16"
}

# Sheet "zh-cn": update Latest Handoff Datetime (col H) for rows 7,10-14
# and set Priority (col E) to "ht" for the same rows
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in 7,10,11,12,13,14) {
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-19 12:20:06"
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
}

# Sheet "de-de": update Latest Handoff Datetime (col H) for rows 7,10-14
# (shares the same underlying text as the Overview sheet's generate date)
# and set Priority (col E) to "ht" for the same rows
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in 7,10,11,12,13,14) {
    $wsDeDe.Cells.Item($r, 8).Value = "2016-08-19 12:20:16"
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
}
